# The "Status" row of the first template table (Status / {{ document.status }})
# is removed, leaving the "{% tr for field in document.fields %}" row as the
# new first row, followed by the field.name / field.value row and the
# "{% tr endfor %}" row. The "_GoBack" bookmark - previously sitting right
# after the "{%p for mandatory in document.mandatoryList %}" paragraph just
# above the second table - is relocated to the very start of that new first
# row's paragraph (this is where Word leaves _GoBack after the edit is made
# directly in that cell).

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$statusRow = $table.Rows.Item(1)
$statusRow.Delete()

$forRow = $table.Rows.Item(1)
$forParagraph = $forRow.Cells.Item(1).Range.Paragraphs.Item(1)
$insertionPoint = $forParagraph.Range.Start

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$newBookmarkRange = $d.Range($insertionPoint, $insertionPoint)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
